$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 205, shifting the existing rows 205-207 down to 206-208.
$ws.Rows.Item(205).Insert()

# Populate the newly inserted row 205 with the new weekly price entry.
$ws.Cells.Item(205, 1).Value = 4
$ws.Cells.Item(205, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(205, 3).Value = "Los Lagos"
$ws.Cells.Item(205, 4).Value = 44628
$ws.Cells.Item(205, 4).NumberFormat = $ws.Cells.Item(206, 4).NumberFormat
$ws.Cells.Item(205, 5).Value = 10
$ws.Cells.Item(205, 6).Value = "Fruta"
$ws.Cells.Item(205, 7).Value = 100109
$ws.Cells.Item(205, 8).Value = "Uva"
$ws.Cells.Item(205, 9).Value = 100109001
$ws.Cells.Item(205, 10).Value = "Uva"
$ws.Cells.Item(205, 11).Value = "Red Globe"
$ws.Cells.Item(205, 12).Value = "Primera"
$ws.Cells.Item(205, 13).Value = 400
$ws.Cells.Item(205, 14).Value = 15000
$ws.Cells.Item(205, 15).Value = 15000
$ws.Cells.Item(205, 16).Value = 15000
$ws.Cells.Item(205, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(205, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(205, 19).Value = 750
$ws.Cells.Item(205, 20).Value = 20
